$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.918.26"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "4.056.07"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.36"
$ws.Range("E5").Value = "  +4.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.96"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "4.050.99"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.699"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.767"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.13"
$ws.Range("E12").Value = "  +13.06%  "
$ws.Range("E13").Value = "  +0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.02"
$ws.Range("E14").Value = "  +2.26%  "
$ws.Range("D15").Value = "4.705.04"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "4.062.42"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.60"
$ws.Range("E17").Value = "  +3.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.88"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("E19").Value = "  +1.14%  "
$ws.Range("E20").Value = "  -0.58%  "
$ws.Range("D21").Value = "72.921.12"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "447.35"
$ws.Range("E22").Value = "  +3.87%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "98.54"
$ws.Range("E23").Value = "  -0.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.56"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.43"
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.36"
$ws.Range("E27").Value = "  +16.47%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.42"
$ws.Range("E28").Value = "  +2.30%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "11.01"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.94"
$ws.Range("E30").Value = "  +1.46%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.46"
$ws.Range("E31").Value = "  +0.88%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.96"
$ws.Range("E32").Value = "  +14.73%  "
$ws.Range("E33").Value = "  +2.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.71"
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "690.33"
$ws.Range("E35").Value = "  +1.25%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "49.25"
$ws.Range("E36").Value = "  +14.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "67.48"
$ws.Range("E37").Value = "  +1.88%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.455"
$ws.Range("E38").Value = "  +6.21%  "
$ws.Range("E39").Value = "  +6.98%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.41"
$ws.Range("E42").Value = "  -2.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.26"
$ws.Range("E43").Value = "  +16.59%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0497"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").Value = "  +4.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.58"
$ws.Range("E49").Value = "  +5.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.33"
$ws.Range("E51").Value = "  -1.87%  "

Write-Output "Applied cryptos update"